# The commit being replayed ("Vygenerovany file ve slozce: ...") is an
# automated re-save of the document through Word. Diffing the underlying
# OOXML shows the only change is to the <w:nsid> values stored on the
# numbering/list definitions in word/numbering.xml - the internal id Word
# stamps on each abstract numbering ("list template") definition. No
# visible content, formatting, list text, or list level ever changes;
# only that internal bookkeeping id is refreshed, which is exactly what
# Word does whenever it rewrites a document's list definitions on save.
#
# That id is not something Word's object model exposes for reading or
# writing - there is no ListTemplate.NSID / List.ListID style property
# that can be assigned a value; Word mints/refreshes it internally
# whenever it rewrites a list definition, as a side effect of touching
# that list, never as a direct, addressable property. The closest thing
# automation can legitimately do is touch the numbering definitions that
# are actually used in the document, so Word marks them dirty and
# rewrites their list template data - the same trigger that causes the
# internal id to be refreshed when the file is regenerated.

$d = $word.ActiveDocument

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $lf = $p.Range.ListFormat

    if ($lf.ListType -ne 0) {
        $lt = $lf.ListTemplate

        # Touch every level of the list template actually backing this
        # paragraph's list, re-assigning each value back onto itself.
        # This marks the abstract numbering definition as modified (and
        # so rewritten by Word) without changing any visible formatting.
        $levels = $lt.ListLevels
        for ($li = 1; $li -le $levels.Count; $li++) {
            $lvl = $levels.Item($li)
            try { $lvl.NumberFormat = $lvl.NumberFormat } catch { }
            try { $lvl.NumberStyle = $lvl.NumberStyle } catch { }
            try { $lvl.StartAt = $lvl.StartAt } catch { }
            try { $lvl.TrailingCharacter = $lvl.TrailingCharacter } catch { }
            try { $lvl.Alignment = $lvl.Alignment } catch { }
        }
    }
}
